$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 ("GS + SB"): the underlying "number of relevant items" (D4) changed
# from 1059 to 1076 (a new informal search using scholar found more hits).
# The dependent ratio/F1-style columns are recomputed from the raw counts:
#   F = C / D
#   G = C / MAX(C)            (unchanged here, C4 did not change)
#   H = harmonic_mean(F, G)
#   I = E / D
#   J = E / MAX(E)            (unchanged here, E4 did not change)
#   K = harmonic_mean(I, J)

$C4 = $ws.Range("C4").Value2
$E4 = $ws.Range("E4").Value2
$G4 = $ws.Range("G4").Value2
$J4 = $ws.Range("J4").Value2

$D4 = 1076
$ws.Range("D4").Value = $D4

$F4 = $C4 / $D4
$I4 = $E4 / $D4
$H4 = (2 * $F4 * $G4) / ($F4 + $G4)
$K4 = (2 * $I4 * $J4) / ($I4 + $J4)

$ws.Range("F4").Value = $F4
$ws.Range("H4").Value = $H4
$ws.Range("I4").Value = $I4
$ws.Range("K4").Value = $K4
